$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name clean-up (drop "(the)"/long-form qualifiers to match ISO short names) ---
$ws.Range("A17").Value  = "Bahamas "
$ws.Range("A27").Value  = "Bolivia"
$ws.Range("A33").Value  = "British Indian Ocean Territory"
$ws.Range("A42").Value  = "Cayman Islands"
$ws.Range("A43").Value  = "Central African Republic"
$ws.Range("A48").Value  = "Cocos Islands"
$ws.Range("A50").Value  = "Comoros"
$ws.Range("A51").Value  = "Democratic Republic of the Congo"
$ws.Range("A52").Value  = "Congo "
$ws.Range("A53").Value  = "Cook Islands"
$ws.Range("A64").Value  = "Dominican Republic "
$ws.Range("A73").Value  = "Falkland Islands"
$ws.Range("A74").Value  = "Faroe Islands "
$ws.Range("A80").Value  = "French Southern Territories "
$ws.Range("A82").Value  = "Gambia "
$ws.Range("A99").Value  = "Holy See"
$ws.Range("A106").Value = "Iran"
$ws.Range("A119").Value = "North Korea"
$ws.Range("A120").Value = "South Korea"
$ws.Range("A123").Value = "Lao People's Democratic Republic"
$ws.Range("A139").Value = "Marshall Islands"
$ws.Range("A145").Value = "Micronesia"
$ws.Range("A146").Value = "Moldova "
$ws.Range("A157").Value = "Netherlands"
$ws.Range("A161").Value = "Niger"
$ws.Range("A165").Value = "Northern Mariana Islands"
$ws.Range("A175").Value = "Philippines"
$ws.Range("A183").Value = "Russia"
$ws.Range("A212").Value = "Sudan"
$ws.Range("A218").Value = "Taiwan"
$ws.Range("A220").Value = "Tanzania"
$ws.Range("A234").Value = "United Arab Emirates"
$ws.Range("A235").Value = "United Kingdom"
$ws.Range("A236").Value = "United States Minor Outlying Islands"
$ws.Range("A237").Value = "United States of America"
$ws.Range("A241").Value = "Venezuela"

# --- Restore explicit 15.75pt row height on the data rows (default changed elsewhere) ---
$ws.Range("1:1000").RowHeight = 15.75

# --- New trailing columns (D:F) sized like the other data columns ---
$ws.Range("D1:F1").ColumnWidth = 11.6
